$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# Start Year: 2020 -> 2050
$ws.Range("B2").Value = 2050

# End Year: 2050 -> 2055
$ws.Range("B3").Value = 2055

# Power_plants_from_year: 2020 -> 2050
$ws.Range("B4").Value = 2050

# Expand the description text for yearly_CO2_prices (row 16 / C16)
$ws.Range("C16").Value = "so far this is only for NL. If False then the price is fixed to the fix_price_year"

# fix_price_year: 2020 -> 2050
$ws.Range("B18").Value = 2050

# fix_demand_to_initial_year: FALSE -> TRUE
$ws.Range("B19").Value = $true

# targetinvestment_per_year active flag: TRUE -> FALSE
$ws.Range("B25").Value = $false

# Row 14 height change (start_tick_dismantling row): 21.5 -> 31.5
$ws.Rows.Item(14).RowHeight = 31.5

# Update the selected cell to match the new active selection C22
$ws.Range("C22").Select()
